# Update "想去人数" (F column) counts in the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Mapping of row number -> new value for column F
$updates = @{
    3  = 109
    4  = 1605
    5  = 615
    7  = 13
    8  = 11435
    9  = 26
    11 = 448
    12 = 356
    13 = 1090
    14 = 794
    15 = 12363
    16 = 13035
    21 = 17
    24 = 104
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
